$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 204793.6
$ws.Range("I15").Value = 204793.6
$ws.Range("K15").Value = 614380.8
$ws.Range("M15").Value = -614211.8
$ws.Range("H38").Value = 220.88889
$ws.Range("I38").Value = 40
$ws.Range("K38").Value = 120
$ws.Range("M38").Value = 252
$ws.Range("H40").Value = 2192.5625
$ws.Range("I40").Value = 2210
$ws.Range("K40").Value = 2210
$ws.Range("M40").Value = -2035
$ws.Range("H42").Value = 100
$ws.Range("I42").Value = 100
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 300
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -70
$ws.Range("N42").Value = $null
$ws.Range("H137").Value = 30304128
$ws.Range("I137").Value = 35715308
$ws.Range("J137").Value = 1520
$ws.Range("K137").Value = 107145924
$ws.Range("L137").Value = 4560
$ws.Range("M137").Value = -107143374
$ws.Range("N137").Value = -9660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6172.6772
$ws.Range("I32").Value = 4354.8765
$ws.Range("J32").Value = 18442.834
$ws.Range("K32").Value = 4354.8765
$ws.Range("L32").Value = 18442.834
$ws.Range("M32").Value = -4067.8765
$ws.Range("N32").Value = -19016.834
$ws.Range("H61").Value = 5751.2
$ws.Range("I61").Value = 5627.76
$ws.Range("J61").Value = 6368.4
$ws.Range("K61").Value = 5627.76
$ws.Range("L61").Value = 6368.4
$ws.Range("M61").Value = -5415.76
$ws.Range("N61").Value = -6792.4
$ws.Range("H74").Value = 6892
$ws.Range("I74").Value = 1200.7142
$ws.Range("J74").Value = 16851.75
$ws.Range("K74").Value = 1200.7142
$ws.Range("L74").Value = 16851.75
$ws.Range("M74").Value = -326.7141999999999
$ws.Range("N74").Value = -18599.75
$ws.Range("H77").Value = 6892
$ws.Range("I77").Value = 1200.7142
$ws.Range("J77").Value = 16851.75
$ws.Range("K77").Value = 6003.571
$ws.Range("L77").Value = 84258.75
$ws.Range("M77").Value = -1635.571
$ws.Range("N77").Value = -92994.75
$ws.Range("H136").Value = 5751.2
$ws.Range("I136").Value = 5627.76
$ws.Range("J136").Value = 6368.4
$ws.Range("K136").Value = 16883.28
$ws.Range("L136").Value = 19105.2
$ws.Range("M136").Value = -14333.28
$ws.Range("N136").Value = -24205.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3478.6572
$ws.Range("I134").Value = 2163.739
$ws.Range("J134").Value = 5998.9165
$ws.Range("K134").Value = 6491.217000000001
$ws.Range("L134").Value = 17996.7495
$ws.Range("M134").Value = -3956.217000000001
$ws.Range("N134").Value = -23066.7495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4034.2712
$ws.Range("I31").Value = 1306.5151
$ws.Range("K31").Value = 1306.5151
$ws.Range("M31").Value = -1011.5151
$ws.Range("H34").Value = 4034.2712
$ws.Range("I34").Value = 1306.5151
$ws.Range("K34").Value = 1306.5151
$ws.Range("M34").Value = -1104.5151
$ws.Range("H54").Value = 7400
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = $null
$ws.Range("H58").Value = 2784.5862
$ws.Range("I58").Value = 1076.1428
$ws.Range("J58").Value = 4379.1333
$ws.Range("K58").Value = 1076.1428
$ws.Range("L58").Value = 4379.1333
$ws.Range("M58").Value = -873.1428000000001
$ws.Range("N58").Value = -4785.1333
$ws.Range("H136").Value = 2784.5862
$ws.Range("I136").Value = 1076.1428
$ws.Range("J136").Value = 4379.1333
$ws.Range("K136").Value = 3228.4284
$ws.Range("L136").Value = 13137.3999
$ws.Range("M136").Value = -678.4284000000002
$ws.Range("N136").Value = -18237.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3529.2727
$ws.Range("I64").Value = 1026
$ws.Range("J64").Value = 4959.7144
$ws.Range("K64").Value = 3078
$ws.Range("L64").Value = 14879.1432
$ws.Range("M64").Value = -2808
$ws.Range("N64").Value = -15419.1432
$ws.Range("H67").Value = 3529.2727
$ws.Range("I67").Value = 1026
$ws.Range("J67").Value = 4959.7144
$ws.Range("K67").Value = 3078
$ws.Range("L67").Value = 14879.1432
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -16751.1432
$ws.Range("H69").Value = 1850
$ws.Range("J69").Value = 2360
$ws.Range("L69").Value = 7080
$ws.Range("N69").Value = -8702
$ws.Range("H72").Value = 1850
$ws.Range("J72").Value = 2360
$ws.Range("L72").Value = 21240
$ws.Range("N72").Value = -29352
$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -47122
$ws.Range("M74").Value = $null
$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -145608
$ws.Range("M77").Value = $null
$ws.Range("H95").Value = 3000
$ws.Range("J95").Value = 3000
$ws.Range("L95").Value = 9000
$ws.Range("N95").Value = -13118
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").Value = $null
$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 836.1
$ws.Range("I22").Value = 606.8333
$ws.Range("J22").Value = 1180
$ws.Range("K22").Value = 606.8333
$ws.Range("L22").Value = 1180
$ws.Range("M22").Value = -311.8333
$ws.Range("N22").Value = -1770
$ws.Range("H27").Value = 836.1
$ws.Range("I27").Value = 606.8333
$ws.Range("J27").Value = 1180
$ws.Range("K27").Value = 606.8333
$ws.Range("L27").Value = 1180
$ws.Range("M27").Value = -499.8333
$ws.Range("N27").Value = -1394
$ws.Range("H46").Value = 1023.875
$ws.Range("I46").Value = 947.5
$ws.Range("J46").Value = 1100.25
$ws.Range("K46").Value = 947.5
$ws.Range("L46").Value = 1100.25
$ws.Range("M46").Value = -759.5
$ws.Range("N46").Value = -1476.25
$ws.Range("H94").Value = 41333.332
$ws.Range("J94").Value = 41333.332
$ws.Range("L94").Value = 41333.332
$ws.Range("N94").Value = -42685.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 50420
$ws.Range("J16").Value = 50420
$ws.Range("L16").Value = 50420
$ws.Range("N16").Value = -50420
